$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.038712620735168
$ws.Range("B1").Value = 1.372332334518433
$ws.Range("C1").Value = 2.193686008453369
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 2.004915475845337
